$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.716369766863604
$ws.Range("C4").Value = 0.772717497952248
$ws.Range("D4").Value = 0.768126880551711
$ws.Range("E4").Value = 0.754584861941685
$ws.Range("F4").Value = 0.87256696213353
$ws.Range("G4").Value = 0.844461583387313
$ws.Range("H4").Value = 0.841389015556442
$ws.Range("I4").Value = 0.65867391008924
$ws.Range("J4").Value = 0.658259163542218
$ws.Range("K4").Value = 0.702279283302204
$ws.Range("L4").Value = 0.770751009217637
$ws.Range("M4").Value = 0.925808443317861
$ws.Range("N4").Value = 0.561959539989434

$ws.Range("D5").Value = 0.438889632437241
$ws.Range("E5").Value = 0.441403577339026
$ws.Range("F5").Value = 0.692723367340905
$ws.Range("H5").Value = 0.511112243128561
$ws.Range("I5").Value = 0.456721110351558
$ws.Range("N5").Value = 0.402399638605333

$ws.Range("D6").Value = 0.433525987240537
$ws.Range("E6").Value = 0.622061067971373
$ws.Range("F6").Value = 0.764034268173467
$ws.Range("G6").Value = 0.627378598315257
$ws.Range("H6").Value = 0.710473338089272
$ws.Range("I6").Value = 0.580152462807448
$ws.Range("J6").Value = 0.529014769509267
$ws.Range("K6").Value = 0.585754346694377
$ws.Range("L6").Value = 0.568451357363088
$ws.Range("M6").Value = 0.880923373592755
$ws.Range("N6").Value = 0.553829259151759
